# Update the two-digit multiplication problems in the worksheet table.
# Each "old=" expression is unique within the document, so a simple
# global Find/Replace (ReplaceAll) per pair is safe and order-independent.

$d = $word.ActiveDocument

$pairs = @(
    @("34×40=", "11×41="),
    @("29×76=", "67×11="),
    @("16×29=", "60×81="),
    @("31×22=", "61×75="),
    @("53×51=", "99×99="),
    @("25×39=", "26×96="),
    @("21×74=", "82×45="),
    @("86×64=", "43×82="),
    @("42×92=", "18×56="),
    @("11×48=", "16×73="),
    @("63×86=", "36×13="),
    @("97×80=", "96×28="),
    @("38×80=", "46×90="),
    @("41×30=", "34×61="),
    @("16×94=", "40×26="),
    @("15×78=", "91×99="),
    @("31×81=", "97×33="),
    @("81×69=", "25×53="),
    @("37×91=", "20×63="),
    @("25×57=", "15×14="),
    @("61×51=", "82×59="),
    @("68×78=", "28×25="),
    @("77×99=", "15×37="),
    @("15×63=", "29×14="),
    @("58×30=", "15×11=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
